$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.042.83"
$ws.Range("E2").Value = "  -4.40%  "

# Row 3
$ws.Range("D3").Value = "3.077.11"
$ws.Range("E3").Value = "  -4.74%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $ws.Range("D3").Style
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.97"
$ws.Range("D5").Style = $ws.Range("D3").Style
$ws.Range("E5").Value = "  -5.98%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.03"
$ws.Range("D6").Style = $ws.Range("D3").Style
$ws.Range("E6").Value = "  -10.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $ws.Range("D3").Style
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("D8").Value = "3.070.70"
$ws.Range("E8").Value = "  -4.63%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("D9").Style = $ws.Range("D3").Style
$ws.Range("E9").Value = "  -4.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("D10").Style = $ws.Range("D3").Style
$ws.Range("E10").Value = "  -4.00%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.24"
$ws.Range("D11").Style = $ws.Range("D3").Style
$ws.Range("E11").Value = "  -11.85%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").Style = $ws.Range("D3").Style
$ws.Range("E12").Value = "  -4.56%  "

# Row 13
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("D13").Style = $ws.Range("D3").Style
$ws.Range("E13").Value = "  -2.68%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.74"
$ws.Range("D14").Style = $ws.Range("D3").Style
$ws.Range("E14").Value = "  -8.16%  "

# Row 15
$ws.Range("D15").Value = "3.544.32"
$ws.Range("E15").Value = "  -5.27%  "

# Row 16
$ws.Range("D16").Value = "62.906.69"
$ws.Range("E16").Value = "  -4.81%  "

# Row 17
$ws.Range("E17").Value = "  -2.95%  "

# Row 18
$ws.Range("D18").Value = "3.075.58"
$ws.Range("E18").Value = "  -4.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.64"
$ws.Range("D19").Style = $ws.Range("D3").Style
$ws.Range("E19").Value = "  -6.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.22"
$ws.Range("D20").Style = $ws.Range("D3").Style
$ws.Range("E20").Value = "  -9.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("D21").Style = $ws.Range("D3").Style
$ws.Range("E21").Value = "  -7.32%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.706"
$ws.Range("D22").Style = $ws.Range("D3").Style
$ws.Range("E22").Value = "  -4.41%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.23"
$ws.Range("D23").Style = $ws.Range("D3").Style
$ws.Range("E23").Value = "  -6.84%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.54"
$ws.Range("D24").Style = $ws.Range("D3").Style
$ws.Range("E24").Value = "  -2.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.13"
$ws.Range("D25").Style = $ws.Range("D3").Style
$ws.Range("E25").Value = "  -9.39%  "

# Row 26
$ws.Range("E26").Value = "  +0.17%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.72"
$ws.Range("D27").Style = $ws.Range("D3").Style
$ws.Range("E27").Value = "  -7.68%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.18"
$ws.Range("D28").Style = $ws.Range("D3").Style
$ws.Range("E28").Value = "  -11.70%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = $ws.Range("D3").Style
$ws.Range("E29").Value = "  -0.53%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("D30").Style = $ws.Range("D3").Style
$ws.Range("E30").Value = "  -14.21%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.09"
$ws.Range("D31").Style = $ws.Range("D3").Style
$ws.Range("E31").Value = "  -4.92%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.10"
$ws.Range("D32").Style = $ws.Range("D3").Style
$ws.Range("E32").Value = "  -6.10%  "

# Row 33
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "59.12"
$ws.Range("D33").Style = $ws.Range("D3").Style
$ws.Range("E33").Value = "  +8.26%  "

# Row 34
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.43"
$ws.Range("D34").Style = $ws.Range("D3").Style
$ws.Range("E34").Value = "  -10.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.01"
$ws.Range("D35").Style = $ws.Range("D3").Style
$ws.Range("E35").Value = "  -4.60%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.15"
$ws.Range("D36").Style = $ws.Range("D3").Style
$ws.Range("E36").Value = "  -7.57%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "473.19"
$ws.Range("D37").Style = $ws.Range("D3").Style
$ws.Range("E37").Value = "  -15.64%  "

# Row 38
$ws.Range("D38").Value = "3.128.93"
$ws.Range("E38").Value = "  -1.48%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0392"
$ws.Range("D39").Style = $ws.Range("D3").Style
$ws.Range("E39").Value = "  -12.88%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0796"
$ws.Range("D40").Style = $ws.Range("D3").Style
$ws.Range("E40").Value = "  -6.60%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.116"
$ws.Range("D41").Style = $ws.Range("D3").Style
$ws.Range("E41").Value = "  -9.61%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.08"
$ws.Range("D42").Style = $ws.Range("D3").Style
$ws.Range("E42").Value = "  -4.98%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("D43").Style = $ws.Range("D3").Style
$ws.Range("E43").Value = "  -10.79%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.252"
$ws.Range("D44").Style = $ws.Range("D3").Style
$ws.Range("E44").Value = "  -10.76%  "

# Row 45
$ws.Range("E45").Value = "  +0.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.03"
$ws.Range("D46").Style = $ws.Range("D3").Style
$ws.Range("E46").Value = "  -11.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.69"
$ws.Range("D47").Style = $ws.Range("D3").Style
$ws.Range("E47").Value = "  -5.87%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.08"
$ws.Range("D48").Style = $ws.Range("D3").Style
$ws.Range("E48").Value = "  -4.66%  "

# Row 49
$ws.Range("E49").Value = "  -3.77%  "

# Row 50
$ws.Range("D50").Value = "0.0₃0514"
$ws.Range("E50").Value = "  -6.62%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.01"
$ws.Range("D51").Style = $ws.Range("D3").Style
$ws.Range("E51").Value = "  -7.85%  "

Write-Host "All edits applied."